# Adds a new "Random Primms" results section (columns AL:AP) to Sheet1,
# mirroring the structure of the existing "Big City Pop" section
# (Voronoi/Primms headers, Pop/Geo sub-headers, 20 data rows, and
# MEAN/STDEV summary formulas).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Section header (row 6) ---
$ws.Range("AL6").Value = "Random Primms"

# --- Sub headers (row 8): Voronoi / Primms ---
$ws.Range("AL8").Value = "Voronoi"
$ws.Range("AO8").Value = "Primms"

# --- Column headers (row 9): Pop / Geo ---
$ws.Range("AL9").Value = "Pop"
$ws.Range("AM9").Value = "Geo"
$ws.Range("AO9").Value = "Pop"
$ws.Range("AP9").Value = "Geo"

# --- Data rows 10-29 ---
$ws.Range("AL10").Value = 1
$ws.Range("AM10").Value = 0.96779999999999999
$ws.Range("AO10").Value = 1
$ws.Range("AP10").Value = 0.87929999999999997

$ws.Range("AL11").Value = 1
$ws.Range("AM11").Value = 0.8649
$ws.Range("AO11").Value = 1
$ws.Range("AP11").Value = 0.78800000000000003

$ws.Range("AL12").Value = 0.625
$ws.Range("AM12").Value = 0.89990000000000003
$ws.Range("AO12").Value = 0.75
$ws.Range("AP12").Value = 0.81320000000000003

$ws.Range("AL13").Value = 0.875
$ws.Range("AM13").Value = 0.9254
$ws.Range("AO13").Value = 0.875
$ws.Range("AP13").Value = 0.84470000000000001

$ws.Range("AL14").Value = 0.75
$ws.Range("AM14").Value = 0.88590000000000002
$ws.Range("AO14").Value = 0.875
$ws.Range("AP14").Value = 0.81889999999999996

$ws.Range("AL15").Value = 0.75
$ws.Range("AM15").Value = 0.9325
$ws.Range("AO15").Value = 0.75
$ws.Range("AP15").Value = 0.84660000000000002

$ws.Range("AL16").Value = 1
$ws.Range("AM16").Value = 0.90469999999999995
$ws.Range("AO16").Value = 1
$ws.Range("AP16").Value = 0.84289999999999998

$ws.Range("AL17").Value = 0.875
$ws.Range("AM17").Value = 0.96360000000000001
$ws.Range("AO17").Value = 0.75
$ws.Range("AP17").Value = 0.89829999999999999

$ws.Range("AL18").Value = 0.875
$ws.Range("AM18").Value = 0.87590000000000001
$ws.Range("AO18").Value = 0.875
$ws.Range("AP18").Value = 0.78979999999999995

$ws.Range("AL19").Value = 0.75
$ws.Range("AM19").Value = 0.86580000000000001
$ws.Range("AO19").Value = 0.875
$ws.Range("AP19").Value = 0.8206

$ws.Range("AL20").Value = 1
$ws.Range("AM20").Value = 0.96640000000000004
$ws.Range("AO20").Value = 1
$ws.Range("AP20").Value = 0.89710000000000001

$ws.Range("AL21").Value = 1
$ws.Range("AM21").Value = 0.96699999999999997
$ws.Range("AO21").Value = 1
$ws.Range("AP21").Value = 0.90290000000000004

$ws.Range("AL22").Value = 0.875
$ws.Range("AM22").Value = 0.81399999999999995
$ws.Range("AO22").Value = 1
$ws.Range("AP22").Value = 0.77329999999999999

$ws.Range("AL23").Value = 1
$ws.Range("AM23").Value = 0.8881
$ws.Range("AO23").Value = 1
$ws.Range("AP23").Value = 0.80559999999999998

$ws.Range("AL24").Value = 1
$ws.Range("AM24").Value = 0.92769999999999997
$ws.Range("AO24").Value = 1
$ws.Range("AP24").Value = 0.85919999999999996

$ws.Range("AL25").Value = 0.875
$ws.Range("AM25").Value = 0.96
$ws.Range("AO25").Value = 0.875
$ws.Range("AP25").Value = 0.877

$ws.Range("AL26").Value = 0.75
$ws.Range("AM26").Value = 0.96440000000000003
$ws.Range("AO26").Value = 0.75
$ws.Range("AP26").Value = 0.88109999999999999

$ws.Range("AL27").Value = 0.75
$ws.Range("AM27").Value = 0.86950000000000005
$ws.Range("AO27").Value = 0.75
$ws.Range("AP27").Value = 0.78259999999999996

$ws.Range("AL28").Value = 1
$ws.Range("AM28").Value = 0.96230000000000004
$ws.Range("AO28").Value = 1
$ws.Range("AP28").Value = 0.87519999999999998

$ws.Range("AL29").Value = 0.875
$ws.Range("AM29").Value = 0.88
$ws.Range("AO29").Value = 0.75
$ws.Range("AP29").Value = 0.7923

# --- Summary rows 30 (MEAN) / 31 (STDEV) ---
$ws.Range("AL30").Formula = "=AVERAGE(AL10:AL29)"
$ws.Range("AM30").Formula = "=AVERAGE(AM10:AM29)"
$ws.Range("AO30").Formula = "=AVERAGE(AO10:AO29)"
$ws.Range("AP30").Formula = "=AVERAGE(AP10:AP29)"

$ws.Range("AL31").Formula = "=STDEV(AL10:AL29)"
$ws.Range("AM31").Formula = "=STDEV(AM10:AM29)"
$ws.Range("AO31").Formula = "=STDEV(AO10:AO29)"
$ws.Range("AP31").Formula = "=STDEV(AP10:AP29)"

# Match the "Good" style (green, used for MEAN/STDEV rows) already applied
# to the other section totals, so the new summary cells look consistent.
# (Skip the AN gap column - the source section pattern does not emit a
# blank styled placeholder cell there, e.g. compare AJ/AL with no AK style.)
$ws.Range("AL30").Style = $ws.Range("AJ30").Style
$ws.Range("AM30").Style = $ws.Range("AJ30").Style
$ws.Range("AO30").Style = $ws.Range("AJ30").Style
$ws.Range("AP30").Style = $ws.Range("AJ30").Style
$ws.Range("AL31").Style = $ws.Range("AJ31").Style
$ws.Range("AM31").Style = $ws.Range("AJ31").Style
$ws.Range("AO31").Style = $ws.Range("AJ31").Style
$ws.Range("AP31").Style = $ws.Range("AJ31").Style

# --- Update selection to match the edited state ---
$ws.Range("AL30").Select()
